$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7466333508491516
$ws.Range("B1").Value = 1.996257543563843
$ws.Range("C1").Value = 3.996880531311035
$ws.Range("D1").Value = 3.536841154098511
$ws.Range("E1").Value = 2.019597291946411
